$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 0.09851898182174161
$ws.Range("C7").Value = 1.540035768533897
$ws.Range("D7").Value = 4.199063518274832
$ws.Range("E7").Value = 2.049161662308475
$ws.Range("F7").Value = 2.074266956696792
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = 0.3078469792503876
$ws.Range("C8").Value = 1.619327053205559
$ws.Range("D8").Value = 4.192265333528787
$ws.Range("E8").Value = 2.047502218198746
$ws.Range("F8").Value = 2.052148847503259
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = 1.149512620956961
$ws.Range("C9").Value = 1.619587756671262
$ws.Range("D9").Value = 4.367097144185886
$ws.Range("E9").Value = 2.08976006856909
$ws.Range("F9").Value = 1.790535902856296
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = 0.8237816712625872
$ws.Range("C10").Value = 1.495310396438245
$ws.Range("D10").Value = 3.378278017657149
$ws.Range("E10").Value = 1.838009253963959
$ws.Range("F10").Value = 1.710156013466232
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = 0.8526019529410073
$ws.Range("C11").Value = 1.762719545188964
$ws.Range("D11").Value = 4.112577759476549
$ws.Range("E11").Value = 2.027949151107233
$ws.Range("F11").Value = 2.05719702183509
$ws.Range("G11").Value = 5
